$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 679.8
$ws.Range("I18").Value = 687.5
$ws.Range("J18").Value = 649
$ws.Range("K18").Value = 687.5
$ws.Range("L18").Value = 649
$ws.Range("M18").Value = -403.5
$ws.Range("N18").Value = -1217
$ws.Range("H40").Value = 5127.8335
$ws.Range("I40").Value = 3744.5
$ws.Range("J40").Value = 7894.5
$ws.Range("K40").Value = 3744.5
$ws.Range("L40").Value = 7894.5
$ws.Range("M40").Value = -3569.5
$ws.Range("N40").Value = -8244.5
$ws.Range("H55").Value = 626.46155
$ws.Range("I55").Value = 37.5
$ws.Range("J55").Value = 733.5454999999999
$ws.Range("K55").Value = 37.5
$ws.Range("L55").Value = 733.5454999999999
$ws.Range("M55").Value = 176.5
$ws.Range("N55").Value = -1161.5455
$ws.Range("H64").Value = 4908.647
$ws.Range("I64").Value = 4521.75
$ws.Range("J64").Value = 5252.5557
$ws.Range("K64").Value = 4521.75
$ws.Range("L64").Value = 5252.5557
$ws.Range("M64").Value = -4273.75
$ws.Range("N64").Value = -5748.5557
$ws.Range("H67").Value = 4908.647
$ws.Range("I67").Value = 4521.75
$ws.Range("J67").Value = 5252.5557
$ws.Range("K67").Value = 4521.75
$ws.Range("L67").Value = 5252.5557
$ws.Range("M67").Value = -3663.75
$ws.Range("N67").Value = -6968.5557
$ws.Range("H74").Value = 12966
$ws.Range("I74").Value = 6957.9165
$ws.Range("K74").Value = 6957.9165
$ws.Range("M74").Value = -6021.9165
$ws.Range("H76").Value = 4962
$ws.Range("J76").Value = 4947.25
$ws.Range("L76").Value = 4947.25
$ws.Range("N76").Value = -5577.25
$ws.Range("H77").Value = 12966
$ws.Range("I77").Value = 6957.9165
$ws.Range("K77").Value = 34789.5825
$ws.Range("M77").Value = -30109.5825
$ws.Range("H79").Value = 4962
$ws.Range("J79").Value = 4947.25
$ws.Range("L79").Value = 4947.25
$ws.Range("N79").Value = -7131.25
$ws.Range("H80").Value = 1042.7142
$ws.Range("I80").Value = 400
$ws.Range("J80").Value = 1149.8334
$ws.Range("K80").Value = 1200
$ws.Range("L80").Value = 3449.5002
$ws.Range("M80").Value = -202
$ws.Range("N80").Value = -5445.5002
$ws.Range("H83").Value = 1042.7142
$ws.Range("I83").Value = 400
$ws.Range("J83").Value = 1149.8334
$ws.Range("K83").Value = 3600
$ws.Range("L83").Value = 10348.5006
$ws.Range("M83").Value = 1392
$ws.Range("N83").Value = -20332.5006
$ws.Range("H88").Value = 6026.316
$ws.Range("I88").Value = 9667
$ws.Range("J88").Value = 5343.6875
$ws.Range("K88").Value = 9667
$ws.Range("L88").Value = 5343.6875
$ws.Range("M88").Value = -9261
$ws.Range("N88").Value = -6155.6875
$ws.Range("H91").Value = 6026.316
$ws.Range("I91").Value = 9667
$ws.Range("J91").Value = 5343.6875
$ws.Range("K91").Value = 9667
$ws.Range("L91").Value = 5343.6875
$ws.Range("M91").Value = -8263
$ws.Range("N91").Value = -8151.6875
$ws.Range("H98").Value = 1294.2667
$ws.Range("I98").Value = 1263.9459
$ws.Range("K98").Value = 1263.9459
$ws.Range("M98").Value = 234.0541000000001
$ws.Range("H105").Value = 24999.5
$ws.Range("J105").Value = 24999.5
$ws.Range("L105").Value = 24999.5
$ws.Range("N105").Value = -31987.5
$ws.Range("H109").Value = 99899.5
$ws.Range("J109").Value = 99899.5
$ws.Range("L109").Value = 99899.5
$ws.Range("N109").Value = -102673.5
$ws.Range("H112").Value = 2434.8
$ws.Range("I112").Value = 1950
$ws.Range("J112").Value = 2460.3157
$ws.Range("K112").Value = 5850
$ws.Range("L112").Value = 7380.9471
$ws.Range("M112").Value = -4742
$ws.Range("N112").Value = -9596.947100000001
$ws.Range("H113").Value = 6313.273
$ws.Range("I113").Value = 5591.5
$ws.Range("K113").Value = 5591.5
$ws.Range("M113").Value = -2337.5
$ws.Range("H116").Value = 13657.546
$ws.Range("J116").Value = 16779.625
$ws.Range("L116").Value = 16779.625
$ws.Range("N116").Value = -23663.625
$ws.Range("H120").Value = 78100
$ws.Range("J120").Value = 78100
$ws.Range("L120").Value = 78100
$ws.Range("N120").Value = -87776
$ws.Range("H122").Value = 1294.2667
$ws.Range("I122").Value = 1263.9459
$ws.Range("K122").Value = 3791.8377
$ws.Range("M122").Value = -1341.8377
$ws.Range("H132").Value = 13426.23
$ws.Range("I132").Value = 15144.637
$ws.Range("K132").Value = 45433.911
$ws.Range("M132").Value = -42903.911
$ws.Range("H135").Value = 1314.9
$ws.Range("I135").Value = 1266.5555
$ws.Range("K135").Value = 11398.9995
$ws.Range("M135").Value = -8863.9995
$ws.Range("H138").Value = 3756
$ws.Range("I138").Value = 3763.4285
$ws.Range("K138").Value = 11290.2855
$ws.Range("M138").Value = -6150.2855
$ws.Range("H141").Value = 3175.5715
$ws.Range("I141").Value = 2886
$ws.Range("K141").Value = 8658
$ws.Range("M141").Value = -3478

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2423.9473
$ws.Range("I102").Value = 1795.9333
$ws.Range("K102").Value = 1795.9333
$ws.Range("M102").Value = -173.9332999999999
$ws.Range("H107").Value = 149333.33
$ws.Range("J107").Value = 149333.33
$ws.Range("L107").Value = 149333.33
$ws.Range("N107").Value = -157013.33
$ws.Range("H109").Value = 26833.666
$ws.Range("J109").Value = 26833.666
$ws.Range("L109").Value = 26833.666
$ws.Range("N109").Value = -29607.666
$ws.Range("H122").Value = 1375.1428
$ws.Range("I122").Value = 1168.16
$ws.Range("J122").Value = 3100
$ws.Range("K122").Value = 3504.48
$ws.Range("L122").Value = 9300
$ws.Range("M122").Value = -1054.48
$ws.Range("N122").Value = -14200
$ws.Range("H132").Value = 60946.94
$ws.Range("I132").Value = 60946.94
$ws.Range("K132").Value = 182840.82
$ws.Range("M132").Value = -180310.82

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 143613.28
$ws.Range("I22").Value = 200658
$ws.Range("K22").Value = 200658
$ws.Range("M22").Value = -200485
$ws.Range("H94").Value = 1601.4333
$ws.Range("I94").Value = 1452.9584
$ws.Range("J94").Value = 2195.3333
$ws.Range("K94").Value = 1452.9584
$ws.Range("L94").Value = 2195.3333
$ws.Range("M94").Value = -1001.9584
$ws.Range("N94").Value = -3097.3333
$ws.Range("H99").Value = 57706.25
$ws.Range("I99").Value = 78582.5
$ws.Range("K99").Value = 78582.5
$ws.Range("M99").Value = -77084.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 39000
$ws.Range("J124").Value = 39000
$ws.Range("L124").Value = 39000
$ws.Range("N124").Value = -43910
$ws.Range("H132").Value = 3330
$ws.Range("I132").Value = 3062
$ws.Range("K132").Value = 9186
$ws.Range("M132").Value = -6656
$ws.Range("H134").Value = 109603.5
$ws.Range("I134").Value = 118669
$ws.Range("K134").Value = 356007
$ws.Range("M134").Value = -353472

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 495.84616
$ws.Range("I121").Value = 384.5
$ws.Range("J121").Value = 545.3333
$ws.Range("K121").Value = 1153.5
$ws.Range("L121").Value = 1635.9999
$ws.Range("M121").Value = 156.5
$ws.Range("N121").Value = -4255.9999
$ws.Range("H125").Value = 4873.75
$ws.Range("I125").Value = 4873.75
$ws.Range("K125").Value = 14621.25
$ws.Range("M125").Value = -9701.25
$ws.Range("H131").Value = 10265.477
$ws.Range("J131").Value = 7630.9287
$ws.Range("L131").Value = 22892.7861
$ws.Range("N131").Value = -32972.7861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1458.8096
$ws.Range("I80").Value = 868.5
$ws.Range("J80").Value = 1822.0769
$ws.Range("K80").Value = 868.5
$ws.Range("L80").Value = 1822.0769
$ws.Range("M80").Value = 129.5
$ws.Range("N80").Value = -3818.0769
$ws.Range("H83").Value = 1458.8096
$ws.Range("I83").Value = 868.5
$ws.Range("J83").Value = 1822.0769
$ws.Range("K83").Value = 4342.5
$ws.Range("L83").Value = 9110.3845
$ws.Range("M83").Value = 649.5
$ws.Range("N83").Value = -19094.3845
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1384.6111
$ws.Range("I16").Value = 1407.2354
$ws.Range("K16").Value = 1407.2354
$ws.Range("M16").Value = -1237.2354
$ws.Range("H55").Value = 89.8
$ws.Range("J55").Value = 84.333336
$ws.Range("L55").Value = 84.333336
$ws.Range("N55").Value = -430.333336
$ws.Range("H93").Value = 2025.2667
$ws.Range("I93").Value = 1618.6666
$ws.Range("J93").Value = 2635.1667
$ws.Range("K93").Value = 1618.6666
$ws.Range("L93").Value = 2635.1667
$ws.Range("M93").Value = -370.6666
$ws.Range("N93").Value = -5131.1667
$ws.Range("H100").Value = 3262.6875
$ws.Range("I100").Value = 2896.4614
$ws.Range("J100").Value = 4849.6665
$ws.Range("K100").Value = 2896.4614
$ws.Range("L100").Value = 4849.6665
$ws.Range("M100").Value = -2355.4614
$ws.Range("N100").Value = -5931.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2222.3125
$ws.Range("I136").Value = 1803.8
$ws.Range("K136").Value = 5411.4
$ws.Range("M136").Value = -2861.4
